$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the special "last row" formatting (thicker bottom border styling) from row 29 onto row 26,
#    since after removing the LUDIN SARAYS employee (rows 27-29), row 26 becomes the new last row
#    of the data table and must carry the bottom-border style.
$ws.Range("B29:J29").Copy() | Out-Null
$ws.Range("B26:J26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2) Update the account-statement data table (rows 16-26) with the new employee/period/value data.
#    Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$ws.Range("C16").Value2 = "1126428231"
$ws.Range("D16").Value2 = "JUAN ESTEBAN ARANGO CADAVID"
$ws.Range("E16").Value2 = "2504"
$ws.Range("F16").Value2 = 58000

$ws.Range("C17").Value2 = "73185307"
$ws.Range("D17").Value2 = "JORGE ALBERTO TORRES CARRASCAL"
$ws.Range("E17").Value2 = "2505"
$ws.Range("F17").Value2 = 22000

$ws.Range("C18").Value2 = "73291014"
$ws.Range("D18").Value2 = "LUIS ALFONSO AVILA ROMERO"
$ws.Range("E18").Value2 = "2505"
$ws.Range("F18").Value2 = 22000

$ws.Range("C19").Value2 = "1126428231"
$ws.Range("D19").Value2 = "JUAN ESTEBAN ARANGO CADAVID"
$ws.Range("E19").Value2 = "2505"
$ws.Range("F19").Value2 = 56000

$ws.Range("C20").Value2 = "1050973070"
$ws.Range("D20").Value2 = "SERGIO ANDRES CERA BOLAÃ?O"
$ws.Range("E20").Value2 = "2505"
$ws.Range("F20").Value2 = 22000

$ws.Range("C21").Value2 = "73185307"
$ws.Range("D21").Value2 = "JORGE ALBERTO TORRES CARRASCAL"
$ws.Range("E21").Value2 = "2506"
$ws.Range("F21").Value2 = 60000

$ws.Range("C22").Value2 = "73291014"
$ws.Range("D22").Value2 = "LUIS ALFONSO AVILA ROMERO"
$ws.Range("E22").Value2 = "2506"
$ws.Range("F22").Value2 = 60000

$ws.Range("C23").Value2 = "1050973070"
$ws.Range("D23").Value2 = "SERGIO ANDRES CERA BOLAÃ?O"
$ws.Range("E23").Value2 = "2506"
$ws.Range("F23").Value2 = 60000

$ws.Range("C24").Value2 = "73185307"
$ws.Range("D24").Value2 = "JORGE ALBERTO TORRES CARRASCAL"
$ws.Range("E24").Value2 = "2507"
$ws.Range("F24").Value2 = 60000

$ws.Range("C25").Value2 = "73291014"
$ws.Range("D25").Value2 = "LUIS ALFONSO AVILA ROMERO"
$ws.Range("E25").Value2 = "2507"
$ws.Range("F25").Value2 = 14000

$ws.Range("C26").Value2 = "1050973070"
$ws.Range("D26").Value2 = "SERGIO ANDRES CERA BOLAÃ?O"
$ws.Range("E26").Value2 = "2507"
$ws.Range("F26").Value2 = 60000

# 3) Remove the now-obsolete LUDIN SARAYS DE LA ROSA GUERRERO rows entirely (this also shifts the
#    signature block rows 34/35 up to 31/32, and shrinks dimension/mergeCells/sharedStrings accordingly).
$ws.Rows("27:29").Delete()

# 4) Update the header summary figures.
$ws.Range("E11").Value2 = 494000
$ws.Range("C13").Value2 = 4
$ws.Range("F13").Value2 = 4

# 5) Column D (Nombre Trabajador) autofit width shrinks now that the longest name
#    ("LUDIN SARAYS DE LA ROSA GUERRERO") was removed.
$ws.Columns("D").ColumnWidth = 34.8
